# "try implement repository pattern"
#
# The workbook is being reshaped so two lookup/junction sheets read better
# as explicit "Films<Something>" link tables:
#   GenreList -> FilmsGenres   (Films<->Genres link table)
#   ListFilms -> FilmsUsers    (Films<->Users link table)
#
# The FilmsGenres sheet also gets its header cells renamed from
# FilmsId/GenresId to IdFilms/IdGenres (matching the IdFilms/IdUser naming
# already used on the FilmsUsers sheet), and becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Rename the two sheets.
$wsGenreList = $wb.Worksheets.Item("GenreList")
$wsGenreList.Name = "FilmsGenres"

$wsListFilms = $wb.Worksheets.Item("ListFilms")
$wsListFilms.Name = "FilmsUsers"

# Update the FilmsGenres header row to the new Id naming convention.
$wsGenreList.Range("A1").Value = "IdFilms"
$wsGenreList.Range("B1").Value = "IdGenres"

# Make FilmsGenres the active sheet/tab with B1 selected.
$wsGenreList.Activate()
$wsGenreList.Range("B1").Select()
